# "Parser verwerkt nu alle vragen die antwoorden hebben"
#
# Sheet2 (the answers grid) used to start its data one column too far to
# the right (col B.."vraag 1" etc.), with a stray "ZSST" value sitting in
# A2. The fix shifts the whole sheet one column to the left (which also
# drops that stray A2 value) and inserts a brand-new question ("vraag 2")
# with two answer rows that were missing from the parser's output.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: shift everything one column to the left -----------------
# Deleting column A shifts B->A, C->B, ... and discards the old A-column
# contents (the stray "ZSST" marker in A2) along with it.
$ws2.Columns.Item(1).Delete()

# --- Sheet2: insert the new "vraag 2" question's two answer rows -----
# After the shift, "vraag 2" now lives in A6 (followed directly by what
# used to be "vraag 2" -> now "vraag 3" two rows later once we insert).
# Make room for its two answers right after row 6.
$ws2.Rows.Item(7).Resize(2).Insert()

# Row 8 first, then row 7, so the new shared-string entries land in the
# same order as in the target workbook (Antwo394u, ajajaja.png, weiorj9).
$ws2.Cells.Item(8, 2).Value = "Antwo394u"
$ws2.Cells.Item(8, 3).Value = "ajajaja.png"
$ws2.Cells.Item(7, 2).Value = "weiorj9"

# "x" marks per answer row (which symptom columns apply).
$ws2.Cells.Item(7, 11).Value = "x"
$ws2.Cells.Item(7, 13).Value = "x"
$ws2.Cells.Item(8, 8).Value  = "x"
$ws2.Cells.Item(8, 15).Value = "x"
$ws2.Cells.Item(8, 18).Value = "x"
$ws2.Cells.Item(8, 20).Value = "x"

# --- Selections --------------------------------------------------------
# Sheet1's saved selection moves to the full header row (A2:XFD2); select
# it on Sheet1 without leaving Sheet2 as the non-active tab.
$ws1.Range("A2:XFD2").Select()
$ws2.Activate()
$ws2.Range("T8").Select()
